$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E1").Value = "field_tag"
$ws.Range("E2").Value = "R008"
$ws.Range("E4").Value = "G058"
$ws.Range("E5").Value = "R013"
$ws.Range("E3").Value = "B028"
$ws.Range("F8").Select()
